# Keytrack_Framework/TestData/Data.xlsx
# "Code pushed of verifying company in side bar and companies drop down."
#
# The author bumped the Gmail "+tag" used for the Login test fixture from
# +44 to +45 (e.g. after re-running/re-recording the login test with a new
# throwaway alias), and the sheet's last active-cell selection was left on
# A4 when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Update the login email used by the test data (A1) -> +45 alias.
$ws.Range("A1").Value = "tahirgeeks+45@gmail.com"

# Leave the sheet's selection/active cell on A4, as captured in the saved view.
[void]$ws.Range("A4").Select()
